$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    # Footers: both the "first page" footer and the "default" footer carry the
    # Pearson logo (wp:docPr/name="image1.png") which needs to become
    # "image2.png".
    foreach ($hf in $sec.Footers) {
        foreach ($sh in $hf.Range.InlineShapes) {
            if ($sh.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $sh.Name = "image2.png"
            }
        }
    }

    # Headers: the "first page" header carries the BTec logo
    # (wp:docPr/name="image2.jpg") which needs to become "image1.jpg".
    foreach ($hf in $sec.Headers) {
        foreach ($sh in $hf.Range.InlineShapes) {
            if ($sh.AlternativeText -eq "BTec_Logo-Orange") {
                $sh.Name = "image1.jpg"
            }
        }
    }
}
